$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (Big Onion): append growing-area place names to "Suitable Areas" (B2)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Mainly in dry zone. Can be cultivated in areas from sea level to 2000 meters.`nMatale and Anuradhapura,`nGalewela, Dambulla, Kimbissa, Maradankadawala,`nPolonnaruwa, Kurunagala, Vavuniya, Mullaitive,`nMannar"

# ---------------------------------------------------------------------------
# Row 3 (Rice): fill in "Suitable Areas" (B3) and "Soil Condition" (F3)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Anuradhapura, Kurunagala, Ampara, Polonnaruwa, Batticaloa, Hambanthota, Monaragala, Trincomalee, Kilinochchi, Badulla, Vavuniya, Mannar, Matale, Mulativu, Puttalam, Ratnapura, Kandy, Matara, Gampaha, Jaffna, Kalutara, Galle, Kegalle, Colombo, NuwaraEliya"
$ws.Range("F3").Value = "The p.h. value between 5.8 to 7.0 is more suitable."

# ---------------------------------------------------------------------------
# Row 4 (Chilli): append growing-area place names to "Suitable Areas" (B4)
# and add harvest-day details to the "Recommended varieties" rich text (G4)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Chilli can be grown from sea level to about 1600 meters high.`nAnuradhapura, Polonnaruwa,`nKurunagala,`nMatale, Puttalam"

$g4 = "K A 02 - Suitable for Dry zone. More suited to growing in the Yala season. Can harvest in 60 to 90 days.`nArunalu - Suitable for Yala season. Can harvest in 60 to 70 days.`nMI HOT - Suitable for Dry and Wet zones. 60 to 70 days for harvest.`nMI GREEN - More suitable for the Northern Province. 60 to 70 days `n"
$ws.Range("G4").Value = $g4
$cellG4 = $ws.Range("G4")

$cellG4.Characters(1, 9).Font.Bold = $true
$cellG4.Characters(1, 9).Font.Name = "Arial"
$cellG4.Characters(1, 9).Font.Size = 11

$cellG4.Characters(10, 96).Font.Bold = $false
$cellG4.Characters(10, 96).Font.Name = "Arial"
$cellG4.Characters(10, 96).Font.Size = 11

$cellG4.Characters(106, 10).Font.Bold = $true
$cellG4.Characters(106, 10).Font.Name = "Arial"
$cellG4.Characters(106, 10).Font.Size = 11

$cellG4.Characters(116, 56).Font.Bold = $false
$cellG4.Characters(116, 56).Font.Name = "Arial"
$cellG4.Characters(116, 56).Font.Size = 11

$cellG4.Characters(172, 9).Font.Bold = $true
$cellG4.Characters(172, 9).Font.Name = "Arial"
$cellG4.Characters(172, 9).Font.Size = 11

$cellG4.Characters(181, 59).Font.Bold = $false
$cellG4.Characters(181, 59).Font.Name = "Arial"
$cellG4.Characters(181, 59).Font.Size = 11

$cellG4.Characters(240, 11).Font.Bold = $true
$cellG4.Characters(240, 11).Font.Name = "Arial"
$cellG4.Characters(240, 11).Font.Size = 11

$cellG4.Characters(251, 56).Font.Bold = $false
$cellG4.Characters(251, 56).Font.Name = "Arial"
$cellG4.Characters(251, 56).Font.Size = 11

# ---------------------------------------------------------------------------
# Row 5 (Maize / Corn): append growing-area place names to "Suitable Areas"
# (B5) and turn the "Recommended varieties" text (G5) into rich text with
# the variety names in bold plus the new harvest-day sentences.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Dry and temperate zones are best suited for commercial cultivation.`nAmpara, Anuradhapura, Polonnaruwa, Kurunagala, Moneragala, Badulla, Matale"

$g5 = "Badhra - Harvest time is between 105-110 days.`nM.I.Mase Hybrid 01 & M.I.Mase Hybrid  - Harvest time is between 105-110 days.`nM.I.Mase Hybrid 03 - Harvest time is between 100-105 days."
$ws.Range("G5").Value = $g5
$cellG5 = $ws.Range("G5")

$cellG5.Characters(1, 7).Font.Bold = $true
$cellG5.Characters(1, 7).Font.Name = "Arial"
$cellG5.Characters(1, 7).Font.Size = 11

$cellG5.Characters(8, 40).Font.Bold = $false
$cellG5.Characters(8, 40).Font.Name = "Arial"
$cellG5.Characters(8, 40).Font.Size = 11

$cellG5.Characters(48, 36).Font.Bold = $true
$cellG5.Characters(48, 36).Font.Name = "Arial"
$cellG5.Characters(48, 36).Font.Size = 11

$cellG5.Characters(84, 42).Font.Bold = $false
$cellG5.Characters(84, 42).Font.Name = "Arial"
$cellG5.Characters(84, 42).Font.Size = 11

$cellG5.Characters(126, 18).Font.Bold = $true
$cellG5.Characters(126, 18).Font.Name = "Arial"
$cellG5.Characters(126, 18).Font.Size = 11

$cellG5.Characters(144, 40).Font.Bold = $false
$cellG5.Characters(144, 40).Font.Name = "Arial"
$cellG5.Characters(144, 40).Font.Size = 11

# ---------------------------------------------------------------------------
# Row 7 (Finger millet): append growing-area place names to "Suitable Areas"
# (B7)
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Southern dry zone`nKurunagala, Anuradhapura, Moneragala, Polonnaruwa, Badulla, Matale"

# ---------------------------------------------------------------------------
# Row 2 grew taller once the extra place names were added - match the
# resulting autofit height recorded in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 254.25

# ---------------------------------------------------------------------------
# Selection / scroll position moved as part of the author's review pass.
# ---------------------------------------------------------------------------
$ws.Range("G5").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
